$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "306.87"
Set-TextValue $ws.Range("E2") "-6.39%"
Set-TextValue $ws.Range("G2") "23"
Set-TextValue $ws.Range("D3") "39.32"
Set-TextValue $ws.Range("E3") "-11.30%"
Set-TextValue $ws.Range("G3") "23"
Set-TextValue $ws.Range("D4") "4.997"
Set-TextValue $ws.Range("E4") "-6.96%"
Set-TextValue $ws.Range("G4") "23"
Set-TextValue $ws.Range("D5") "0.07745"
Set-TextValue $ws.Range("E5") "-7.41%"
Set-TextValue $ws.Range("G5") "23"
Set-TextValue $ws.Range("D6") "4.285"
Set-TextValue $ws.Range("E6") "-3.43%"
Set-TextValue $ws.Range("G6") "23"
Set-TextValue $ws.Range("D7") "1.594"
Set-TextValue $ws.Range("E7") "-18.02%"
Set-TextValue $ws.Range("G7") "23"
Set-TextValue $ws.Range("D8") "0.9184"
Set-TextValue $ws.Range("E8") "-5.52%"
Set-TextValue $ws.Range("G8") "23"
Set-TextValue $ws.Range("D9") "0.09927"
Set-TextValue $ws.Range("E9") "-13.02%"
Set-TextValue $ws.Range("G9") "23"
Set-TextValue $ws.Range("D10") "0.1730"
Set-TextValue $ws.Range("E10") "-9.26%"
Set-TextValue $ws.Range("G10") "23"
Set-TextValue $ws.Range("D11") "0.08914"
Set-TextValue $ws.Range("E11") "-7.36%"
Set-TextValue $ws.Range("G11") "23"
Set-TextValue $ws.Range("D12") "0.04385"
Set-TextValue $ws.Range("E12") "-5.24%"
Set-TextValue $ws.Range("G12") "23"
Set-TextValue $ws.Range("D13") "7.037"
Set-TextValue $ws.Range("E13") "-16.68%"
Set-TextValue $ws.Range("G13") "23"
Set-TextValue $ws.Range("D14") "0.1057"
Set-TextValue $ws.Range("E14") "-0.33%"
Set-TextValue $ws.Range("G14") "23"
Set-TextValue $ws.Range("D15") "0.001247"
Set-TextValue $ws.Range("E15") "-3.98%"
Set-TextValue $ws.Range("G15") "23"
Set-TextValue $ws.Range("D16") "0.005650"
Set-TextValue $ws.Range("E16") "-4.65%"
Set-TextValue $ws.Range("G16") "23"
Set-TextValue $ws.Range("D17") "3.365"
Set-TextValue $ws.Range("E17") "-0.03%"
Set-TextValue $ws.Range("G17") "23"
Set-TextValue $ws.Range("D18") "2.588"
Set-TextValue $ws.Range("E18") "2.20%"
Set-TextValue $ws.Range("G18") "23"
Set-TextValue $ws.Range("D19") "0.3367"
Set-TextValue $ws.Range("E19") "0.29%"
Set-TextValue $ws.Range("G19") "23"
Set-TextValue $ws.Range("D20") "0.1364"
Set-TextValue $ws.Range("E20") "0.81%"
Set-TextValue $ws.Range("G20") "23"
Set-TextValue $ws.Range("D21") "0.2776"
Set-TextValue $ws.Range("E21") "1.88%"
Set-TextValue $ws.Range("G21") "23"
Set-TextValue $ws.Range("D22") "0.04127"
Set-TextValue $ws.Range("E22") "-1.25%"
Set-TextValue $ws.Range("G22") "23"
Set-TextValue $ws.Range("D23") "0.001204"
Set-TextValue $ws.Range("E23") "-3.06%"
Set-TextValue $ws.Range("G23") "23"
Set-TextValue $ws.Range("D24") "0.004083"
Set-TextValue $ws.Range("E24") "-7.70%"
Set-TextValue $ws.Range("G24") "23"
Set-TextValue $ws.Range("D25") "0.0001226"
Set-TextValue $ws.Range("E25") "-5.85%"
Set-TextValue $ws.Range("G25") "23"
Set-TextValue $ws.Range("D26") "0.0002992"
Set-TextValue $ws.Range("E26") "0.24%"
Set-TextValue $ws.Range("G26") "23"
Set-TextValue $ws.Range("G27") "23"
Set-TextValue $ws.Range("G28") "23"
Set-TextValue $ws.Range("G29") "23"
Set-TextValue $ws.Range("G30") "23"
Set-TextValue $ws.Range("G31") "23"
Set-TextValue $ws.Range("G32") "23"
Set-TextValue $ws.Range("G33") "23"
Set-TextValue $ws.Range("G34") "23"
Set-TextValue $ws.Range("G35") "23"
Set-TextValue $ws.Range("G36") "23"
Set-TextValue $ws.Range("G37") "23"
Set-TextValue $ws.Range("D38") "0.02360"
Set-TextValue $ws.Range("E38") "-12.52%"
Set-TextValue $ws.Range("G38") "23"
Set-TextValue $ws.Range("E39") "-8.81%"
Set-TextValue $ws.Range("G39") "23"
Set-TextValue $ws.Range("D40") "0.007994"
Set-TextValue $ws.Range("E40") "3.29%"
Set-TextValue $ws.Range("G40") "23"
Set-TextValue $ws.Range("D41") "0.1327"
Set-TextValue $ws.Range("E41") "-5.91%"
Set-TextValue $ws.Range("G41") "23"
Set-TextValue $ws.Range("D42") "0.007382"
Set-TextValue $ws.Range("E42") "0.79%"
Set-TextValue $ws.Range("G42") "23"
Set-TextValue $ws.Range("D43") "0.002003"
Set-TextValue $ws.Range("E43") "-5.64%"
Set-TextValue $ws.Range("G43") "23"
Set-TextValue $ws.Range("D44") "0.007265"
Set-TextValue $ws.Range("E44") "-16.71%"
Set-TextValue $ws.Range("G44") "23"
Set-TextValue $ws.Range("D45") "0.3323"
Set-TextValue $ws.Range("E45") "-5.21%"
Set-TextValue $ws.Range("G45") "23"
Set-TextValue $ws.Range("D46") "0.00006705"
Set-TextValue $ws.Range("E46") "-2.94%"
Set-TextValue $ws.Range("G46") "23"
Set-TextValue $ws.Range("D47") "0.00000000754"
Set-TextValue $ws.Range("E47") "0.31%"
Set-TextValue $ws.Range("G47") "23"
Set-TextValue $ws.Range("D48") "0.003328"
Set-TextValue $ws.Range("E48") "-5.00%"
Set-TextValue $ws.Range("G48") "23"
Set-TextValue $ws.Range("D49") "0.004118"
Set-TextValue $ws.Range("E49") "16.43%"
Set-TextValue $ws.Range("G49") "23"
Set-TextValue $ws.Range("D50") "0.00002110"
Set-TextValue $ws.Range("E50") "0.31%"
Set-TextValue $ws.Range("G50") "23"
Set-TextValue $ws.Range("D51") "0.0002010"
Set-TextValue $ws.Range("E51") "0.31%"
Set-TextValue $ws.Range("G51") "23"

Write-Host "Applied all cell updates"
